$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.921.47'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '2.306.43'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '306.38'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.07'
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.72'
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0797'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '18.54'
$ws.Range('E12').Value = '  +4.68%  '
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.78'
$ws.Range('E14').Value = '  -1.69%  '
$ws.Range('D15').Value = '2.666.57'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '2.307.45'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.784'
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('D18').Value = '42.862.58'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.02'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').Value = '0.0₃0903'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('E21').Value = '  -1.50%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.57'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.85'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '25.47'
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '167.19'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.10'
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '33.20'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.77'
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '17.96'
$ws.Range('E35').Value = '  -1.52%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.00'
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('E41').Value = '  -0.98%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.72'
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('D43').Value = '2.007.27'
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0280'
$ws.Range('E44').Value = '  -2.54%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '18.32'
$ws.Range('E45').Value = '  +3.69%  '
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  -3.86%  '
$ws.Range('E48').Value = '  -2.55%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.94'
$ws.Range('E49').Value = '  +7.65%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '53.88'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').Value = '2.534.75'
$ws.Range('E51').Value = '  -0.19%  '
